$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string / cell text for C15 (append statistics note)
$ws.Range("C15").Value = "Päivämäärän hallinta toimii nyt, alotettu draw.io dokumentaatio, statseja Nba komponenttiin"

# Replace B15's literal value with a formula summing the extra work items
$ws.Range("B15").Formula = "=80+34+16+30+60"

# Row 15 gets taller (wrapped text) to accommodate the longer note
$ws.Rows.Item(15).RowHeight = 30

# Move the active selection from H15 to M10
$ws.Range("M10").Select()
